$wb = $excel.ActiveWorkbook

# The two new truck sheets are duplicates (same layout/styles/conditional
# formatting) of the existing "Bus_Makhulu_r" sheet, renamed and with the
# "Instance" cell (H3) updated to the new name, inserted at the end of the
# workbook (this matches how the author created them in Excel: duplicate the
# rear-differential template sheet for each new rigid-axle truck).
$template = $wb.Worksheets.Item("Bus_Makhulu_r")

$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$truckA2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$truckA2.Name = "Truck_Amandla_A2"
$truckA2.Range("H3").Value = "Gear1DShafts1D_Truck_Amandla_A2"

$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$truckA3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$truckA3.Name = "Truck_Amandla_A3"
$truckA3.Range("H3").Value = "Gear1DShafts1D_Truck_Amandla_A3"

# The newly-added, last sheet becomes the active / selected tab.
$truckA3.Activate()
